$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for rows 2-11 from
# 2023-10-08 (serial 45207) to 2023-10-09 (serial 45208).
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45208
}
